$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: batch_013
$ws.Range("A14").Value = "batch_013"
$ws.Range("B14").Value = "y"
$ws.Range("C14").Value = "批量操作语句13执行"
$ws.Range("D14").Value = "batchsql"
$ws.Range("E14").Value = "SingleTable"
$ws.Range("G14").Value = "batch013"
$ws.Range("I14").Value = "batch_sql_013"
$ws.Range("J14").Value = "select * from `$batch013"
$ws.Range("K14").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_013.csv"
$ws.Range("N14").Value = "csv_containsAll"

# Row 15: batch_014
$ws.Range("A15").Value = "batch_014"
$ws.Range("B15").Value = "y"
$ws.Range("C15").Value = "批量操作语句14执行"
$ws.Range("D15").Value = "batchsql"
$ws.Range("E15").Value = "SingleTable"
$ws.Range("G15").Value = "batch014"
$ws.Range("I15").Value = "batch_sql_014"
$ws.Range("J15").Value = "select * from `$batch014"
$ws.Range("K15").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_014.csv"
$ws.Range("N15").Value = "csv_containsAll"

# Row 16: batch_015
$ws.Range("A16").Value = "batch_015"
$ws.Range("B16").Value = "y"
$ws.Range("C16").Value = "批量操作语句15执行"
$ws.Range("D16").Value = "batchsql"
$ws.Range("E16").Value = "SingleTable"
$ws.Range("G16").Value = "batch015"
$ws.Range("I16").Value = "batch_sql_015"
$ws.Range("J16").Value = "select id,name,age from `$batch015"
$ws.Range("K16").Value = "src/test/resources/io.dingodb.test/testdata/cases/batchsql/expectedresult/batch_015.csv"
$ws.Range("N16").Value = "csv_containsAll"

# Match source formatting: K column cells use "fill" horizontal alignment.
$ws.Range("K14").HorizontalAlignment = 5
$ws.Range("K15").HorizontalAlignment = 5
$ws.Range("K16").HorizontalAlignment = 5

$ws.Range("G24").Select()
